$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "clicking on boxview cheat pertection"
$ws.Range("C5").Value = "user can only click on the first row at the start and then work their way down the grid only allowed to click on theat curretn row they are on"
$ws.Range("D5").Value = "check pertechion works"
$ws.Range("E5").Value = "PASS"

$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("C6").Select()
